$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "Recommandations": rows 2-38 get refreshed BRVM/stock figures;
# rows 39-46 (now obsolete) are cleared so the used range shrinks to G38.
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Recommandations")

$data1 = @(
    ,@("NEI-CEDA CI", 0, 3, 2835, 990, "🟡 Observer", "➖ Neutre")
    ,@("BRVM - SERVICES PUBLICS", 0, 6, 2471.88, 103.71, "🟡 Observer", "➖ Neutre")
    ,@("BRVM - AUTRES SECTEURS", 0, 3, 1883.97, 657.38, "🟡 Observer", "➖ Neutre")
    ,@("BRVM - DISTRIBUTION", 0, 3, 1475.25, 492.94, "🟡 Observer", "➖ Neutre")
    ,@("BRVM - TRANSPORT", 0, 3, 1053.75, 354.92, "🟡 Observer", "➖ Neutre")
    ,@("BRVM - AGRICULTURE", 0, 3, 986.31, 330.24, "🟡 Observer", "➖ Neutre")
    ,@("BRVM - INDUSTRIE  (**)", 0, 3, 769.56, 255.18, "🟡 Observer", "➖ Neutre")
    ,@("BRVM-PRINCIPAL  (**)", 0, 3, 648.02, 215.48, "🟡 Observer", "➖ Neutre")
    ,@("BRVM - CONSOMMATION DE BASE  (**)", 0, 3, 640.29, 212.61, "🟡 Observer", "➖ Neutre")
    ,@("BRVM - CONSOMMATION DISCRETIONNAIRE", 0, 3, 510.72, 172.35, "🟡 Observer", "➖ Neutre")
    ,@("BRVM - FINANCES", 0, 3, 435.66, 145, "🟡 Observer", "➖ Neutre")
    ,@("BRVM - SERVICES FINANCIERS", 0, 3, 428.17, 142.51, "🟡 Observer", "➖ Neutre")
    ,@("BRVM-PRESTIGE", 0, 3, 419.07, 139.83, "🟡 Observer", "➖ Neutre")
    ,@("BRVM - INDUSTRIELS", 0, 3, 388.94, 131.67, "🟡 Observer", "➖ Neutre")
    ,@("BRVM - ENERGIE", 0, 3, 329.09, 109.64, "🟡 Observer", "➖ Neutre")
    ,@("BRVM - TELECOMMUNICATIONS", 0, 3, 280.76, 93.58, "🟡 Observer", "➖ Neutre")
    ,@("SETAO CI (STAC)", 3, 0, 22.06, 7.26, "🟢 Achat", "✅ Renforcer")
    ,@("ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)", 2, 0, 14.72, 7.25, "🟡 Observer", "➖ Neutre")
    ,@("NEI-CEDA CI (NEIC)", 2, 0, 9.76, 4.76, "🟡 Observer", "➖ Neutre")
    ,@("LOTERIE NATIONALE DU BENIN (LNBB)", 1, 0, 5.13, 5.13, "🟡 Observer", "➖ Neutre")
    ,@("ORAGROUP TOGO (ORGT)", 1, 0, 2.6, 2.6, "🟡 Observer", "➖ Neutre")
    ,@("TRACTAFRIC MOTORS CI (PRSC)", 1, 0, 1.6, 1.6, "🟡 Observer", "➖ Neutre")
    ,@("SICABLE CI (CABC)", 1, 1, 1.32, 7.5, "🟡 Observer", "👀 À surveiller")
    ,@("TOTAL", 0, 3, 0, 0, "🟡 Observer", "➖ Neutre")
    ,@("CFAO MOTORS CI (CFAC)", 1, 1, -0.24, 3.21, "🟡 Observer", "👀 À surveiller")
    ,@("SUCRIVOIRE (SCRC)", 1, 1, -0.36, 2.78, "🟡 Observer", "👀 À surveiller")
    ,@("SAPH CI (SPHC)", 1, 1, -0.73, 4.93, "🟡 Observer", "👀 À surveiller")
    ,@("BANK OF AFRICA NG (BOAN)", 1, 1, -1.68, -7.45, "🟡 Observer", "➖ Neutre")
    ,@("BANK OF AFRICA ML (BOAM)", 0, 1, -2.68, -2.68, "🟡 Observer", "➖ Neutre")
    ,@("SOLIBRA CI (SLBC)", 0, 1, -2.74, -2.74, "🟡 Observer", "➖ Neutre")
    ,@("NSIA BANQUE COTE D'IVOIRE (NSBC)", 0, 1, -3.51, -3.51, "🟡 Observer", "➖ Neutre")
    ,@("BANK OF AFRICA BF (BOABF)", 0, 1, -3.55, -3.55, "🟡 Observer", "➖ Neutre")
    ,@("CIE CI (CIEC)", 0, 1, -4.09, -4.09, "🟡 Observer", "➖ Neutre")
    ,@("SODE CI (SDCC)", 0, 1, -4.92, -4.92, "🟡 Observer", "➖ Neutre")
    ,@("SICOR CI (SICC)", 0, 1, -7.42, -7.42, "🟡 Observer", "➖ Neutre")
    ,@("UNILEVER CI (UNLC)", 0, 1, -7.5, -7.5, "🟡 Observer", "➖ Neutre")
    ,@("VIVO ENERGY CI (SHEC)", 0, 2, -12.11, -6.74, "🟡 Observer", "➖ Neutre")
)

for ($i = 0; $i -lt $data1.Length; $i++) {
    $r = 2 + $i
    $row = $data1[$i]
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Cells.Item($r, 7).Value = $row[6]
}

# Drop the 8 trailing rows that no longer exist (46 -> 38 rows).
$ws1.Range("A39:G46").ClearContents()

# ----------------------------------------------------------------------
# Sheet "Top_YTD": refreshed YTD progression figures. Rows 2-7 keep their
# sector/stock labels; rows 8-11 are relabeled to the new top sectors.
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Top_YTD")

$ws2.Cells.Item(2, 2).Value = 477836.07
$ws2.Cells.Item(3, 2).Value = 113805
$ws2.Cells.Item(4, 2).Value = 38427.42
$ws2.Cells.Item(5, 2).Value = 20619.72
$ws2.Cells.Item(6, 2).Value = 9088.04
$ws2.Cells.Item(7, 2).Value = 7782.22

$ws2.Cells.Item(8, 1).Value = "BRVM - INDUSTRIE  (**)"
$ws2.Cells.Item(8, 2).Value = 4431.55

$ws2.Cells.Item(9, 1).Value = "BRVM-PRINCIPAL  (**)"
$ws2.Cells.Item(9, 2).Value = 3055.63

$ws2.Cells.Item(10, 1).Value = "BRVM - CONSOMMATION DE BASE  (**)"
$ws2.Cells.Item(10, 2).Value = 2979.05

$ws2.Cells.Item(11, 1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws2.Cells.Item(11, 2).Value = 1873.41
